$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C32").Value = "heyyy"
